# Weekly fruit/vegetable data refresh: rows 4-8 are rotated so that the
# former "Región de O'Higgins / granel" entries (previously dated 44309,
# rows 7-8) now lead as rows 4-5, and the former "Provincia del Elquí /
# empedrada" entries (previously dated 44285, rows 4-6) follow as rows 6-8.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4 (was: Especial/Elqui/empedrada dated 44285 -> now: Primera/O'Higgins/granel dated 44309)
$ws.Range("D4").Value = 44309
$ws.Range("L4").Value = "Primera"
$ws.Range("Q4").Value = "$/caja 15 kilos granel"
$ws.Range("R4").Value = "Región de O'Higgins"

# Row 5 (was: Primera/Elqui/empedrada dated 44285 -> now: Segunda/O'Higgins/granel dated 44309)
$ws.Range("D5").Value = 44309
$ws.Range("L5").Value = "Segunda"
$ws.Range("M5").Value = 70
$ws.Range("Q5").Value = "$/caja 15 kilos granel"
$ws.Range("R5").Value = "Región de O'Higgins"

# Row 6 (was: Segunda/Elqui/empedrada -> now: Especial/Elqui/empedrada, same date/source but new quality & volume/price)
$ws.Range("L6").Value = "Especial"
$ws.Range("M6").Value = 40
$ws.Range("N6").Value = 18000
$ws.Range("O6").Value = 18000
$ws.Range("P6").Value = 18000
$ws.Range("S6").Value = 1200

# Row 7 (was: Primera/O'Higgins/granel dated 44309 -> now: Primera/Elqui/empedrada dated 44285)
$ws.Range("D7").Value = 44285
$ws.Range("M7").Value = 90
$ws.Range("N7").Value = 15000
$ws.Range("O7").Value = 15000
$ws.Range("P7").Value = 15000
$ws.Range("Q7").Value = "$/caja 15 kilos empedrada"
$ws.Range("R7").Value = "Provincia del Elquí"
$ws.Range("S7").Value = 1000

# Row 8 (was: Segunda/O'Higgins/granel dated 44309 -> now: Segunda/Elqui/empedrada dated 44285)
$ws.Range("D8").Value = 44285
$ws.Range("M8").Value = 75
$ws.Range("N8").Value = 12000
$ws.Range("O8").Value = 12000
$ws.Range("P8").Value = 12000
$ws.Range("Q8").Value = "$/caja 15 kilos empedrada"
$ws.Range("R8").Value = "Provincia del Elquí"
$ws.Range("S8").Value = 800
